$wb = $excel.ActiveWorkbook

# --- Sheet "Suivi de projet" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Suivi de projet")

# Row 7 - "Plan du Reseau": assigned to Julien, actual start/end + work filled in
$ws1.Range("D7").Value = "Julien"
$ws1.Range("I7").Value = 43507
$ws1.Range("J7").Value = 43508
$ws1.Range("L7").Value = 3

# Row 8 - "Packet Tracer": actual end date filled in
$ws1.Range("J8").Value = 43489

# Row 9 - "Installation Windows Server": estimated/actual start moved, actual end filled in
$ws1.Range("E9").Value = 43542
$ws1.Range("I9").Value = 43542
$ws1.Range("J9").Value = 43542

# Row 10 - "Partage Ressources": estimated start switched from placeholder text to a real date
$ws1.Range("E10").Value = 43546

# Row 11 - "Construction Reseau VM": assigned to Julien, placeholder start text updated
$ws1.Range("D11").Value = "Julien"
$ws1.Range("E11").Value = "/"

# Header legend: third name changed from Sofian Roger to Jeremy Masse
$ws1.Range("G1").Value = "JEREMY MASSE"

# Row 12 - "Accès Reseau": assigned to Julien, note added
$ws1.Range("D12").Value = "Julien"
$ws1.Range("O12").Value = "Problème connexion APACHE"

# Row 13 - "WordPress + Dossier": assigned to Julien, estimated start filled in
$ws1.Range("D13").Value = "Julien"
$ws1.Range("E13").Value = 43487

$ws1.Activate()
$ws1.Range("O13").Select()

# --- Sheet "Installation" ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("Installation")

$ws2.Range("C6").Value = "Julien"
$ws2.Range("C7").Value = "Julien"
$ws2.Range("C8").Value = "Julien"
$ws2.Range("C9").Value = "Julien"
$ws2.Range("C10").Value = "Julien"
$ws2.Range("C11").Value = "Julien"

$ws2.Activate()
$ws2.Range("C12").Select()

$ws1.Activate()
